$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(29, 8).Value = 534  # H29: 474.18182 -> 534
$ws.Cells.Item(29, 9).Value = 38.25  # I29: 42 -> 38.25
$ws.Cells.Item(29, 10).Value = 4500  # J29: 1626.6666 -> 4500
$ws.Cells.Item(29, 11).Value = 114.75  # K29: 126 -> 114.75
$ws.Cells.Item(29, 12).Value = 13500  # L29: 4879.9998 -> 13500
$ws.Cells.Item(29, 13).Value = 166.25  # M29: 155 -> 166.25
$ws.Cells.Item(29, 14).Value = -14062  # N29: -5441.9998 -> -14062

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 8026.3335  # H32: 12329.6 -> 8026.3335
$ws.Cells.Item(32, 9).Value = 10091.167  # I32: 14187 -> 10091.167
$ws.Cells.Item(32, 10).Value = 3896.6667  # J32: 4900 -> 3896.6667
$ws.Cells.Item(32, 11).Value = 10091.167  # K32: 14187 -> 10091.167
$ws.Cells.Item(32, 12).Value = 3896.6667  # L32: 4900 -> 3896.6667
$ws.Cells.Item(32, 13).Value = -9765.166999999999  # M32: -13861 -> -9765.166999999999
$ws.Cells.Item(32, 14).Value = -4548.6667  # N32: -5552 -> -4548.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 6198.8887  # H43: 6761.25 -> 6198.8887
$ws.Cells.Item(43, 10).Value = 4466.3335  # J43: 5019.6 -> 4466.3335
$ws.Cells.Item(43, 12).Value = 4466.3335  # L43: 5019.6 -> 4466.3335
$ws.Cells.Item(43, 14).Value = -4604.3335  # N43: -5157.6 -> -4604.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 6160  # H62: 5248.5557 -> 6160
$ws.Cells.Item(62, 9).Value = 4499  # I62: 4997.5 -> 4499
$ws.Cells.Item(62, 10).Value = 6436.8335  # J62: 5320.2856 -> 6436.8335
$ws.Cells.Item(62, 11).Value = 4499  # K62: 4997.5 -> 4499
$ws.Cells.Item(62, 12).Value = 6436.8335  # L62: 5320.2856 -> 6436.8335
$ws.Cells.Item(62, 13).Value = -3875  # M62: -4373.5 -> -3875
$ws.Cells.Item(62, 14).Value = -7684.8335  # N62: -6568.2856 -> -7684.8335

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(65, 8).Value = 6160  # H65: 5248.5557 -> 6160
$ws.Cells.Item(65, 9).Value = 4499  # I65: 4997.5 -> 4499
$ws.Cells.Item(65, 10).Value = 6436.8335  # J65: 5320.2856 -> 6436.8335
$ws.Cells.Item(65, 11).Value = 22495  # K65: 24987.5 -> 22495
$ws.Cells.Item(65, 12).Value = 32184.1675  # L65: 26601.428 -> 32184.1675
$ws.Cells.Item(65, 13).Value = -19375  # M65: -21867.5 -> -19375
$ws.Cells.Item(65, 14).Value = -38424.1675  # N65: -32841.428 -> -38424.1675

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 106240.4  # H86: 106289.8 -> 106240.4
$ws.Cells.Item(86, 9).Value = 15500  # I86: 10666.333 -> 15500
$ws.Cells.Item(86, 10).Value = 128925.5  # J86: 147271.28 -> 128925.5
$ws.Cells.Item(86, 11).Value = 15500  # K86: 10666.333 -> 15500
$ws.Cells.Item(86, 12).Value = 128925.5  # L86: 147271.28 -> 128925.5
$ws.Cells.Item(86, 13).Value = -14377  # M86: -9543.333000000001 -> -14377
$ws.Cells.Item(86, 14).Value = -131171.5  # N86: -149517.28 -> -131171.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 106240.4  # H89: 106289.8 -> 106240.4
$ws.Cells.Item(89, 9).Value = 15500  # I89: 10666.333 -> 15500
$ws.Cells.Item(89, 10).Value = 128925.5  # J89: 147271.28 -> 128925.5
$ws.Cells.Item(89, 11).Value = 77500  # K89: 53331.665 -> 77500
$ws.Cells.Item(89, 12).Value = 644627.5  # L89: 736356.4 -> 644627.5
$ws.Cells.Item(89, 13).Value = -71884  # M89: -47715.665 -> -71884
$ws.Cells.Item(89, 14).Value = -655859.5  # N89: -747588.4 -> -655859.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 23307.941  # H106: 22207.389 -> 23307.941
$ws.Cells.Item(106, 9).Value = 23307.941  # I106: 22207.389 -> 23307.941
$ws.Cells.Item(106, 11).Value = 23307.941  # K106: 22207.389 -> 23307.941
$ws.Cells.Item(106, 13).Value = -22676.941  # M106: -21576.389 -> -22676.941

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 2373.739  # H132: 2542.762 -> 2373.739
$ws.Cells.Item(132, 9).Value = 1300.35  # I132: 1337.2632 -> 1300.35
$ws.Cells.Item(132, 10).Value = 9529.666999999999  # J132: 13995 -> 9529.666999999999
$ws.Cells.Item(132, 11).Value = 3901.05  # K132: 4011.7896 -> 3901.05
$ws.Cells.Item(132, 12).Value = 28589.001  # L132: 41985 -> 28589.001
$ws.Cells.Item(132, 13).Value = -1371.05  # M132: -1481.7896 -> -1371.05
$ws.Cells.Item(132, 14).Value = -33649.001  # N132: -47045 -> -33649.001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2138.2727  # H138: 2051.8 -> 2138.2727
$ws.Cells.Item(138, 9).Value = 1382.7778  # I138: 1420 -> 1382.7778
$ws.Cells.Item(138, 10).Value = 3044.8667  # J138: 2648.5 -> 3044.8667
$ws.Cells.Item(138, 11).Value = 4148.3334  # K138: 4260 -> 4148.3334
$ws.Cells.Item(138, 12).Value = 9134.6001  # L138: 7945.5 -> 9134.6001
$ws.Cells.Item(138, 13).Value = 991.6665999999996  # M138: 880 -> 991.6665999999996
$ws.Cells.Item(138, 14).Value = -19414.6001  # N138: -18225.5 -> -19414.6001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(53, 8).Value = 0  # H53: 4000 -> 0
$ws.Cells.Item(53, 9).Value = 0  # I53: 4000 -> 0
$ws.Cells.Item(53, 11).Value = 0  # K53: 4000 -> 0
$ws.Cells.Item(53, 13).ClearContents()  # M53: -3318 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6425.4287  # H61: 6426.143 -> 6425.4287
$ws.Cells.Item(61, 9).Value = 6663.0835  # I61: 6906.5454 -> 6663.0835
$ws.Cells.Item(61, 10).Value = 4999.5  # J61: 4664.6665 -> 4999.5
$ws.Cells.Item(61, 11).Value = 6663.0835  # K61: 6906.5454 -> 6663.0835
$ws.Cells.Item(61, 12).Value = 4999.5  # L61: 4664.6665 -> 4999.5
$ws.Cells.Item(61, 13).Value = -6451.0835  # M61: -6694.5454 -> -6451.0835
$ws.Cells.Item(61, 14).Value = -5423.5  # N61: -5088.6665 -> -5423.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(62, 8).Value = 0  # H62: 80000 -> 0
$ws.Cells.Item(62, 10).Value = 0  # J62: 80000 -> 0
$ws.Cells.Item(62, 12).ClearContents()  # L62: 80000 -> (removed)
$ws.Cells.Item(62, 14).Value = 0  # N62: -81248 -> 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(65, 8).Value = 0  # H65: 80000 -> 0
$ws.Cells.Item(65, 10).Value = 0  # J65: 80000 -> 0
$ws.Cells.Item(65, 12).ClearContents()  # L65: 240000 -> (removed)
$ws.Cells.Item(65, 14).Value = 0  # N65: -246240 -> 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(119, 8).Value = 41566  # H119: 39849 -> 41566
$ws.Cells.Item(119, 10).Value = 41566  # J119: 39849 -> 41566
$ws.Cells.Item(119, 12).Value = 41566  # L119: 39849 -> 41566
$ws.Cells.Item(119, 14).Value = -51242  # N119: -49525 -> -51242

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 6425.4287  # H136: 6426.143 -> 6425.4287
$ws.Cells.Item(136, 9).Value = 6663.0835  # I136: 6906.5454 -> 6663.0835
$ws.Cells.Item(136, 10).Value = 4999.5  # J136: 4664.6665 -> 4999.5
$ws.Cells.Item(136, 11).Value = 19989.2505  # K136: 20719.6362 -> 19989.2505
$ws.Cells.Item(136, 12).Value = 14998.5  # L136: 13993.9995 -> 14998.5
$ws.Cells.Item(136, 13).Value = -17439.2505  # M136: -18169.6362 -> -17439.2505
$ws.Cells.Item(136, 14).Value = -20098.5  # N136: -19093.9995 -> -20098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(15, 8).Value = 8500  # H15: 8550 -> 8500
$ws.Cells.Item(15, 9).Value = 8500  # I15: 0 -> 8500
$ws.Cells.Item(15, 10).Value = 0  # J15: 8550 -> 0
$ws.Cells.Item(15, 11).Value = 8500  # K15: 0 -> 8500
$ws.Cells.Item(15, 12).ClearContents()  # L15: 8550 -> (removed)
$ws.Cells.Item(15, 13).Value = -8273  # M15: None -> -8273
$ws.Cells.Item(15, 14).Value = 0  # N15: -9004 -> 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 8388.25  # H105: 9529.429 -> 8388.25
$ws.Cells.Item(105, 9).Value = 6729.5713  # I105: 8681.4 -> 6729.5713
$ws.Cells.Item(105, 10).Value = 19999  # J105: 11649.5 -> 19999
$ws.Cells.Item(105, 11).Value = 6729.5713  # K105: 8681.4 -> 6729.5713
$ws.Cells.Item(105, 12).Value = 19999  # L105: 11649.5 -> 19999
$ws.Cells.Item(105, 13).Value = -4982.5713  # M105: -6934.4 -> -4982.5713
$ws.Cells.Item(105, 14).Value = -23493  # N105: -15143.5 -> -23493

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 1084746.9  # H6: 1924775.6 -> 1084746.9
$ws.Cells.Item(6, 9).Value = 300296.7  # I6: 1364371.6 -> 300296.7
$ws.Cells.Item(6, 11).Value = 300296.7  # K6: 1364371.6 -> 300296.7
$ws.Cells.Item(6, 13).Value = -300183.7  # M6: -1364258.6 -> -300183.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 0  # H86: 4999 -> 0
$ws.Cells.Item(86, 9).Value = 0  # I86: 4999 -> 0
$ws.Cells.Item(86, 11).Value = 0  # K86: 4999 -> 0
$ws.Cells.Item(86, 13).ClearContents()  # M86: -3876 -> (removed)

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 0  # H89: 4999 -> 0
$ws.Cells.Item(89, 9).Value = 0  # I89: 4999 -> 0
$ws.Cells.Item(89, 11).Value = 0  # K89: 24995 -> 0
$ws.Cells.Item(89, 13).ClearContents()  # M89: -19379 -> (removed)

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 8).Value = 9750  # H57: 9900 -> 9750
$ws.Cells.Item(57, 10).Value = 9499  # J57: 9874.75 -> 9499
$ws.Cells.Item(57, 12).Value = 28497  # L57: 29624.25 -> 28497
$ws.Cells.Item(57, 14).Value = -29615  # N57: -30742.25 -> -29615

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 125428.875  # H107: 111506.11 -> 125428.875
$ws.Cells.Item(107, 10).Value = 167072.83  # J107: 143223 -> 167072.83
$ws.Cells.Item(107, 12).Value = 501218.49  # L107: 429669 -> 501218.49
$ws.Cells.Item(107, 14).Value = -505058.49  # N107: -433509 -> -505058.49

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(129, 8).Value = 588.3333  # H129: 754.2857 -> 588.3333
$ws.Cells.Item(129, 9).Value = 532.5  # I129: 776 -> 532.5
$ws.Cells.Item(129, 11).Value = 1597.5  # K129: 2328 -> 1597.5
$ws.Cells.Item(129, 13).Value = 3402.5  # M129: 2672 -> 3402.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(137, 8).Value = 9999  # H137: 4799.5 -> 9999
$ws.Cells.Item(137, 9).Value = 9999  # I137: 5000 -> 9999
$ws.Cells.Item(137, 10).Value = 0  # J137: 4599 -> 0
$ws.Cells.Item(137, 11).Value = 29997  # K137: 15000 -> 29997
$ws.Cells.Item(137, 12).Value = 0  # L137: 13797 -> 0
$ws.Cells.Item(137, 13).ClearContents()  # M137: -9900 -> (removed)
$ws.Cells.Item(137, 14).Value = -24897  # N137: -23997 -> -24897

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(35, 8).Value = 2947480  # H35: 2360984 -> 2947480
$ws.Cells.Item(35, 10).Value = 1750000  # J35: 882500 -> 1750000
$ws.Cells.Item(35, 12).Value = 1750000  # L35: 882500 -> 1750000
$ws.Cells.Item(35, 14).Value = -1750596  # N35: -883096 -> -1750596

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 9450.429  # H113: 15916.667 -> 9450.429
$ws.Cells.Item(113, 9).Value = 4701.5  # I113: 7000 -> 4701.5
$ws.Cells.Item(113, 10).Value = 11350  # J113: 17700 -> 11350
$ws.Cells.Item(113, 11).Value = 4701.5  # K113: 7000 -> 4701.5
$ws.Cells.Item(113, 12).Value = 11350  # L113: 17700 -> 11350
$ws.Cells.Item(113, 13).Value = -2531.5  # M113: -4830 -> -2531.5
$ws.Cells.Item(113, 14).Value = -15690  # N113: -22040 -> -15690

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1457.5  # H22: 1529.5834 -> 1457.5
$ws.Cells.Item(22, 10).Value = 1491.8889  # J22: 1625.2858 -> 1491.8889
$ws.Cells.Item(22, 12).Value = 1491.8889  # L22: 1625.2858 -> 1491.8889
$ws.Cells.Item(22, 14).Value = -2081.8889  # N22: -2215.2858 -> -2081.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 1457.5  # H27: 1529.5834 -> 1457.5
$ws.Cells.Item(27, 10).Value = 1491.8889  # J27: 1625.2858 -> 1491.8889
$ws.Cells.Item(27, 12).Value = 1491.8889  # L27: 1625.2858 -> 1491.8889
$ws.Cells.Item(27, 14).Value = -1705.8889  # N27: -1839.2858 -> -1705.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 5000  # H40: 0 -> 5000
$ws.Cells.Item(40, 9).Value = 5000  # I40: 0 -> 5000
$ws.Cells.Item(40, 11).Value = 5000  # K40: 0 -> 5000
$ws.Cells.Item(40, 13).Value = -4864  # M40: None -> -4864

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3664.2856  # H46: 3210.5881 -> 3664.2856
$ws.Cells.Item(46, 9).Value = 2585.7144  # I46: 2398.75 -> 2585.7144
$ws.Cells.Item(46, 10).Value = 4742.857  # J46: 3932.2222 -> 4742.857
$ws.Cells.Item(46, 11).Value = 2585.7144  # K46: 2398.75 -> 2585.7144
$ws.Cells.Item(46, 12).Value = 4742.857  # L46: 3932.2222 -> 4742.857
$ws.Cells.Item(46, 13).Value = -2397.7144  # M46: -2210.75 -> -2397.7144
$ws.Cells.Item(46, 14).Value = -5118.857  # N46: -4308.2222 -> -5118.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48, 8).Value = 0  # H48: 5000 -> 0
$ws.Cells.Item(48, 9).Value = 0  # I48: 5000 -> 0
$ws.Cells.Item(48, 11).Value = 0  # K48: 5000 -> 0
$ws.Cells.Item(48, 13).ClearContents()  # M48: -4339 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 5062.778  # H122: 5006.9644 -> 5062.778
$ws.Cells.Item(122, 9).Value = 3599.8572  # I122: 3587.375 -> 3599.8572
$ws.Cells.Item(122, 11).Value = 10799.5716  # K122: 10762.125 -> 10799.5716
$ws.Cells.Item(122, 13).Value = -8349.571599999999  # M122: -8312.125 -> -8349.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 6690.3  # H136: 6364.4287 -> 6690.3
$ws.Cells.Item(136, 9).Value = 6816.4443  # I136: 6504 -> 6816.4443
$ws.Cells.Item(136, 10).Value = 5555  # J136: 5527 -> 5555
$ws.Cells.Item(136, 11).Value = 20449.3329  # K136: 19512 -> 20449.3329
$ws.Cells.Item(136, 12).Value = 16665  # L136: 16581 -> 16665
$ws.Cells.Item(136, 13).Value = -17899.3329  # M136: -16962 -> -17899.3329
$ws.Cells.Item(136, 14).Value = -21765  # N136: -21681 -> -21765

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(64, 8).Value = 66500  # H64: 66000 -> 66500
$ws.Cells.Item(64, 10).Value = 66500  # J64: 66000 -> 66500
$ws.Cells.Item(64, 12).Value = 66500  # L64: 66000 -> 66500
$ws.Cells.Item(64, 14).Value = -66996  # N64: -66496 -> -66996

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(67, 8).Value = 66500  # H67: 66000 -> 66500
$ws.Cells.Item(67, 10).Value = 66500  # J67: 66000 -> 66500
$ws.Cells.Item(67, 12).Value = 66500  # L67: 66000 -> 66500
$ws.Cells.Item(67, 14).Value = -68216  # N67: -67716 -> -68216

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2361.6875  # H122: 2515.6667 -> 2361.6875
$ws.Cells.Item(122, 10).Value = 4000  # J122: 3915.8333 -> 4000
$ws.Cells.Item(122, 12).Value = 12000  # L122: 11747.4999 -> 12000
$ws.Cells.Item(122, 14).Value = -16900  # N122: -16647.4999 -> -16900
